$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 'Chanson du matin, chanson du soir Chanson d''oiseau, chanson de cigale Chanson à texte fait sourire Nous fait danser, nous fait réfléchir Chanson restée au fond d''un tiroir Chanson classée dans le hit-parade Chanson populaire, chanson de merde Chanson de merde Chanson de merde Chanson festive, mélancolique Chanson grossière, chanson poétique À toi toutes les TS du trottoir L''hymne des piliers de comptoir Chanson du râteau, du roulage de pelle Chanson paillarde, chanson de Noël Chanson d''amour, chanson de merde Chanson de lutte, chanson politique Contre le régime nous pousse à l''exil Chanson fait plier l''oppresseur Fait grincer les dents du dictateur Chanson victime de la censure Subit la foudre et la torture Majeur en l''air, chanson t''emmerde'
$ws.Cells.Item(3, 3).Value = 'C''est mon come-back au quartier Putain, ce con m''avait tant manqué Je suis heureux de le retrouver Je me sens enfin libéré C''est mon come-back au quartier Heureux les gars de vous retrouver Quelque chose en moi s''est rallumé On se met bien, c''est ma tournée Elle m''a dit, prends tes affaires, retourne chez ta mère C''st sur le palier, prends tout t reprends la mer C''est pas l''hôtel, ou quoi une garçonnière Je n''suis pas ta mégère, encore moins ta cuisinière T''as passé l''âge de jouer les débauchés Maintenant ça rentre au bercail parce que Monsieur est fauché Tu regrettes tes écarts mais là c''est trop tard Ramasse tes caleçons, va rejoindre le club des connards Au quartier Oui je rentre au quartier Je suis heureux de le retrouver Je me sens enfin libéré C''est mon come-back au quartier Heureux les gars de vous retrouver Quelque chose en moi s''est rallumé On se met bien, c''est ma tournée Ils m''ont dit, eh Toko Blaze, tu es libérable Range ton paquetage et surtout rends-nous toutes les armes Coup d''pied au cul, pour toi c''est hors de la base Mais qu''est-ce qui t''a pris de tirer sur ton caporal chef Le chef veut devenir le grand chef Des milliers de p''tits chefs prêts à donner leurs fesses Le cul se lèche pour une place au mess Bon vent jeunesse C''est mon come-back au quartier Putain, ce con m''avait tant manqué Je suis heureux de le retrouver Je me sens enfin libéré C''est mon come-back au quartier Heureux les gars de vous retrouver Quelque chose en moi s''est rallumé On se met bien, c''est ma tournée Le coach m''a dit, on te conservera pas C''est le grand choc, l''aventure s''arrête là Grand espoir du quartier, désillusion d''un jeune crack Le meilleur au milieu, en défense, en attaque Gros sacrifice, gros sur le cur de quitter ma région Masquer son chagrin la douleur, surmonter la pression Trêve ou à terre, chercher la compétition J''ai perdu la passion C''est mon come-back au quartier C''est mon come-back au quartier C''est mon come-back au quartier Putain, ce con m''avait tant manqué Je suis heureux de le retrouver Je me sens enfin libéré C''est mon come-back au quartier Heureux les gars de vous retrouver Quelque chose en moi s''est rallumé On se met bien, c''est ma tournée'
$ws.Cells.Item(4, 3).Value = 'On est naïf, on est gosse, angoisse quand on s''endort, la mater veut nous couper les locks C''est la coupe des voyous, le look des marginaux, mais surtout des bab'' et des tox Pour toi deux heures de colle, ça sent la beuh dans l''école, on t''accuse parce que t''as des dreadlocks Toujours coursé par les cops, t''échappes pas au contrôle, sois agile et rusé comme el fox On est peace, on est love, on est rough, on est tough, trop bagarreurs, secs comme un os Teigneux comme un rude boy avec cocktails Molotov, l''opposé de l''image du dreadlocks Comme un vent de révolte, pas d''lyrics à l''eau d''rose, punky reggae party dans l''juke-box On a bleuté les notes, africanisé la France, viens danser sur le son des dreadlocks Dreadlocks what you wanna do, what you gonna do Dreadlocks what you wanna do, what you gonna do Dreadlocks what you wanna do, what you gonna do, what you gonna say Dreadlocks L''État tire passe les cancres à la boxe Au volant de l''écomobile Ça fait vibrer le mur, rendez-vous en festoche On contrôle la zone nous comme des boss Faire rouler la danse pour nous c''est fastoche Envoie le c''est pas du baloche On sort le grand jeu gars c''est comme au cinoche Pour mes poulets, pour mes poulettes qui se roulent des grosse galoches On est sage, on est apaisé on s''endort, un demi-siècle à porter des dreadlocks Préjugés sur l''vieux fou, amateur de psylo, poète dresseur d''éléphants roses Mon éducation ma culture mes diplômes, on s''en moque, que d''a priori sur mes dreadlocks Chercheur ingénieur pilote avocat médecin ou prof, un jour place Beauvau y aura un dreadlocks Dreadlocks what you wanna do, what you gonna do Dreadlocks what you wanna do, what you gonna do Dreadlocks what you wanna do, what you gonna do, what you gonna say Dreadlocks'
$ws.Cells.Item(5, 3).Value = 'J''écris donc je suis, eux suivent la tendance La tendance sent le rance, loin du sens de la vie Érection du bras droit, j''vis pas dans cette France Si ce jour est à eux, j''fais banco sur la nuit Petit-fils d''immigrés, la mémoire dans le sens C''est bien la milice française qui a mit mes deux grands-pères dans les camps Famille italienne logée au quartier Saint-Jean C''est ça ! Raflée au beau-matin, Opération Sultan J''aurais pu caressr la haine, cultiver la vengance Mais je n''ai que poésie et mélodie pour engeance Cela dit, chacun son spot mais tiens la distance Ou tu risques de voir c''que ta mère te montrait dans les films J''fais pas de danses à la con, à deux doigts de l''abîme Ou me faire taper sur l''épaule par le Président de la ? C''est ça ! Dans mon monde, on argumente pas avec toi J''ai ni le temps, ni l''énergie consacré à faire ça Ce quotidien on le capte pas si on vit à Versailles Les grands ont tous cinq doigts mais ils parlent avec trois Anti putes à buzz et allergiques aux bruits de battes J''fricotte avec ces prods qui donnent pull-up et ?shot On fait un good job, pas un tour de chauffe Un coup de chatte, d''un MC en Gucci avec un goût de chiottes Tchaaa, 1.3, Planète Mars Bouay, C.V Ouais x3 Classique comme une paire de Nike Attends j''récupère le mic J''suis full Frozen, Cali ? Si ils tirent, c''est sûr que je vais leur mettre le strike strike La mentale et du toupet Ça part comme de la machin pas coupé T''as vu là j''suis dans lappel groupé Merci les khos, ceux qui ont pas douté Tu prends pipi dans le sauna J''ai investi plus qu''une Daytona Que ça veut eh-eh, elle sait c''qu''on a La j''suis en feu-eu comme Fofana Un béret, une canne Je les choque comme Biggie Je ne fais que du sale, du flagrant délit Tu test, on te bousille, on te fuck, on te bannis De côté j''ai du lourd, demande à Solda Guirri Biggie On caillasse, la police, les balances sont punis On les mets dans la cave et on les tabassent, c''est horrible On t''arrange, tu allonges, la mallette, la valise Si tu payes pas à temps, dans ta voiture on te carbonise Fais le malin, fais le malade, fais le canard, fais le zonard Fais la mala, c''est mes amis qui te feront dans l''Audi noire Tu nous verras, tu nous paieras, tu mourras, on laissera pas de preuves On t''allumera, on est sans foi ni loi, on est des bandits J''représente tous les chiens de la casse, tous les chiens de la casse J''voulais remplir mes murs de disques d''or, mais j''ai pas la place Et je voulais, faire le tour de l''Atlas Ils font pas le quart d''mon taff, et ils veulent ma place Si ça rafale, ça débite, c''qui te tue, c''est les potes Et ceux qui jalousent, c''est les putes, c''est quoi cette putain d''époque ? Je les épates, coup de batte, on suit pas, on suce pas Goodbye, goodnight, c''est quoi cette putain d''époque ? J''porte un calibre à la taille comme une montre au poignet J''te parle pas de ménage quand j''te dis J''nettoyais On a surfait sur la vague où certains se sont noyés La promenade devient une cage, trop de crânes se font broyés J''fume la frappe de Nador, balafré comme Albator Si j''avoue jamais mes tords, c''est parce que je n''en ai pas Avant de parler, agis d''abord, on te gifle si t''es pas d''accord Si on tenvoie la corde, c''est pour t''attacher avec J''ai des armes si t''as des pecs, les douilles peuvent faire tomber des armoires Je côtoie la mort, c''est pour ça que j''ai le regard noir T''as plein d''histoires, j''vais finir par t''appeler Père Castor Zehma ?, mets toi à l''affût et casse-toi Je vois rouge sous les ditas J''arrive en T-Max kité Rafale dans l''habitacle La Terre, plusieurs l''ont quittés On arrive fonsdés, on va niquer le game, choquer, choquer le rap français Ta bitch, la faire danser, la baiser, la faire pioncer Finir et r''commencer, fais pas le fou, on va te tarter Sont personne dans leurs quartiers J''peux pas mourir comme Biggie Ou comme Tupac Amaru Oh mon Andalouse Tous ils sont jaloux Sans pression comme Vinicius, y''a trop de gadjis, faut un bus Tout le monde me dit T''abuses, j''leur répond Il me faut un Russe Elle veut du Venetta, Loro Piana, la poupetta Mais bon, j''suis trop pété, j''m''endors sur son gros terma Ça critique vers la droite, ça critique vers la gauche Ils ont qu''à s''mettre des doigts, j''pense qu''à remplir ma sacoche Marseillais comme la L2, la A7, l''A55 On est 152, imagine si tu nos feintes haha La il me faut 400 plaques Sans prendre le RS4 sans plaque On arrive en équipe le projet on est 404 Sur le booster on est quatre sans casque On les braque sans masque C''est Marseille bébé T''es bonne je te fais un bébé J''ai des armes scellés par les keufs mais j''ai payé mon avocat avec de la végé Sur la vie de ma mère Et de but en blanc, poto on les bute en blanc Mon poto c''est le 13 Toujours en équipe sanglants On les braque sans cagoule sans gants La putain de sa mère C''est Marseille bébé T''es bonne je te fais un bébé J''ai des armes scellés par les keufs mais j''ai payé mon avocat avec de la végé Hé ! Un jour où l''autre on t''aura, on te rafale, te laisse parterre comme dans Gomorra On mangeait le pain sec, pas de sauce Amora, Marseille c''est la mafia, 1.3 la Camorra Avec toi on coupe les ponts c''est définitif, je rêve pas de gros culs, je rêve d''Infinity Elle a des gros seins, cache ça dans ton soutif, elle est trop matrixé un peu comme Trinity Tu veux ma vie ? Prend tout ce qu''il y a avec. J''en veux pas je te la donne Retour d''Hollande, t''es revenu avec 50 grammes, on dirait que t''as fais passer la tonne Oh comique Quand t''as rien, pas de pécule, pas de parloir, c''est pas de leau la taule C''est pas de leau Quand t''es dedans, tout le monde tourne le dos et ta gadji fait sa folle Capitale de la drogue, drogue, drogue, yes sir Du Frozen, que ça smoke, smoke, smoke, yes sir C''est le new 13''Organisé, essaie pas de nous rivaliser Et au pire si on vient en grosse équipe, c''est comme ça c''est la cité RS6 plaqué Pologne mais le bolide est allemand Dans tous les dièses ça charbonne, on se fait plaisir après le Khedma Tu nous fais l''ancien, oui zehma celui qui brasse Ça t''as vu au Kyliad avec Assia de Nasdas 1.3, 5.7, Python, pas de Cobra, tout dans la tête fiston, pas de gros bras Tu te colles de partout, tu donnes ton cul, R''key? sur R''key?, tu es dans de beaux draps Securitas dans le Porsche, pas de menta, ça troue des corps chez nous, pas de Van Damme Ça pense à ta mort, y''a pas de renta, et quand t''es parterre y''a pas de mandat C''est le 13, organe électrique, on va te niquer toi et tes piques DP cagoule levée c''est pratique, la prod je la lève, je la cric Pour tout les cafards postés dans le bloc, ça fout des mini-Uzi dans le froc Pour tout les quartiers où y''a la frappe, nous c''est le 13, on est dans le truc Dans le caleçon, des boules, dans le chargeur, des balles Y''a heja déboule, y''a les keufs détalle Marseille danger Malaga, Rotter'', Alger, harraga J''détaille, j''détalle, j''ai la dalle, chargé, Arai, bécane go Air Max, T-Max, j''démarre fort, pull crache nique tous tes morts gros Organisés comme des pros, t''es solo, on est trop On a fait TP dans la zone la, la, la On a tout niqué pour la somme la, la, la On a esquivé les matons la, la, la 13''Organisé, personne va nous canaliser Lame aiguisée, je souris comme Adebisi, tu vas souffrir selon les organes visés Fais attention à ton intonation, t''es une salope y''a pas d''autres appellations T''es bon qu''à parler, t''as jamais fait l''action, ils sont pas doués en rap, ils sont doués en fellation Tester ? Vaut mieux éviter. Sans hésiter, je suis venu viser les points vitaux J''ai vu tomber la neige en été, j''vois défiler les chiennes, à croire j''suis un véto Appelle le dieci-siete 17, on arrive, commence à t''inquiéter Toi tu parles de la rue comme si tu y étais, la nuit j''entends chuchoter les métaux En GTS Panamera, dans le ghetto la paix se fait rare Nique aujourd''hui, demain on verra, j''ai coffré les albums dans les Tera Ça passe à lacte, ça colle, ça appuie, tout le monde sur tes côtés quand t''as le terrain Un peu de love, ça fait le beau-temps la pluie, fais pas confiance, fais pas l''erreur J''suis pas d''ces bâtards la qui courent, plein de phrases à 20K quand c''est bourré Mental en acier, j''réfléchis trois fois, sur un coup de nerf j''allais mourir Allumé, allumé, allumé, c''est la loi de la jungle, tu te fais allumer Y''a trop de traitres et de folles en effet, Marseille oui tu le sais gros c''est animé Faut le cash pour le Riyadh à ''Kech, c''est pour ça que l''instru j''l''a saccage Ils sont moches comme des mecs qui ont vendus la mèche, nous on va la mettre même si ils sont à 10 dans la cage Un 10 pour la route qui nous paie, un 10 en bas de vos écrans y''en a pas On fait bande à part avec l''argent qu''on a pas, et le talent qu''on ment pas Dos à dos comme Kappa, seuls on va soulever la Copa Intégrer le quota, la proc'' nous met coupable, avant que les coups partent C''est que du business, la cagoule est noire comme les gants, tellement je bombarde, je vais me faire flasher Ne donnes pas trop tes plans aux gens, en équipe ils viendront pour t''éclater C''est que du business, j''te fais pas confiance, teste moi tu vas voir j''fais pas que chanter Le regard est noir comme les vitres teintées, 280 sur la route, toute Alicante Vroom, vroom, tu entends que le moteur, j''éclate tellement fort qu''elle veut pas me lâcher Boom, boom, juste en bas de la tess, celui qui sort est revenu pour manger La miss elle a pas trouvée le prince charmant, c''est mes rats de la tess qui vont t''aborder De Lyon à Marseille c''est la même rengaine, dégaine noire dans le noir pour vesqui les condés D''or et de platine, le cur est black, dans la selha, les pus les blocs J''finis l''taff, j''ai repris ça le matin, gros j''vais pas t''apprendre tous les codes C''est paire de Gore et c''est paire d''Oakley, quand tu parles faut en faire autant Quand faut test ça va pas te rater, y''a des hmm qui partent pas à blanc Faut pas quitter le ghetto, ya un tas de rouleaux sous cello'' Nique, mes sapes sont sous scellé quand je recompte du vert, du yellow Mes gros sait gérer les pagailles, y''a des pes-pom, des libéraux Faut pas jouer les fêlés, faut pas jouer les héros Ça va te mettre la coupe au carré, dans la ville du ballon rond Moi c''est Lil Jah Prod, mais comme Marley j''te fais un p''tit pont Virage DP, Vélodrome chaud comme volcan, si j''fais ça c''est pour la mama Si tu veux faire le bonhomme, tu prends le souffle d''Hiroshima J''arrive à l''abordage avec l''équipage, dans l''instru, dans le découpage Lil Jah, l''Ovni et Sauzer on débarque, sur une scène moi j''ai pas le trac T''es pas concentré, on te rattrape, je crois en Dieu et son karma Rebelle comme Ché Guevara, pour les p''tits de l''Estaque à Saint-Jean la Puenta Bien le bonjour c''est la came qu''on vend Ici, même B-John? sonne normal, j''suis busy, j''sors un gun comme normal J''les bousille, en Prada sur la chaussée, un coup de fil et bientôt un p''tit vient te schlasser Comme à l''époque des Croisades, tu te manges un KO quand la lune est en croissant T''avais jamais vu ça depuis Jacque Mes'', un pétard qui provient de Bucarest Un négro connu pour des faits d''armes, des salopes qui nous font les Fidel Castro Et y''en a des tonnes et des tonnes, toi le jour où tu portes une Daytona je doute Je les vois tous cavaler quand ça détonne, j''les ai vu avaler, qui ça étonne ? C''est Marseille ta grand-mère, ça fait du sale, ça maîtrise pas la grammaire J''suis dans soucis, j''sers les camés, sale pute c''est d''la fraîche la con d''sa ? ta mère Wesh bébé, j''viens d''faire six chiffres c''est le bingo Marseille bébé, on parle pas avec les gringos La plus belle, j''la vois, j''suis khapta sur l''allée J''mets la première et andale, j''suis dans les îles sous Bob Marley C''est Marseille en be-bom, madame veut tester le pe-pom Elle veut LV, Gucci, mais j''vais finir dans le deux-temps On arrive en travers dans le 4Matic, ils ont du mal à respirer comme des asthmatiques On change pas l''équipe toujours la même tactique, posté sur le terrain avec un joint de Static Y''a pas de balles factices, que ça pactise, pas sûr que tu t''en tireras Ce n''est pas dit que le plus fort durera, ça te retrouvera même caché dans le Jura On est cramés comme la plaque de ?, pour la plata gros le petit t''attache Ils sont déterminés, tu leur parles pas d''âge, nan poto y''a plus de places pour le partage Nan ce n''est pas la même pour tous, on sait que pour les notres ils ont plus d''estime Un jour où l''autre poto tout se paye, comme ceux qui font subir à la Palestine La rue c''est pas C8, ça veut tout contrôler comme le G8 Tu fais rien pour ta mère mais tu claques tout pour ta pute qui consomme plus qu''un 4.2L V8 J''sais que les hommes de parole ça court pas les rues, j''te parle c''est sincère je ne joue pas de rôle On peut le faire poto même si ça paraît rude, t''as merdé fils de pute n''inverse pas les rôles C''est pas Marseille, c''est le ralliement, c''est les gars d''ici qui flambent C''est les bâtiments, c''est les grandes baraques, c''est le p''tit bar qui fait l''angle Le diable qui joue les anges, ça fleurit sous les lampes C''est trop inné, on arrose les pieds de la Skunk à l''eau du robinet Je s''rais le gominé, en Stone l''hiver, on peut tous monter, ça va pas se livrer Ça prolifère, rouge dans les yeux, sous le conifère, que des gros tickets J''baise des travailleuses, pris dans l''ivresse, j''ai le feu sous les doigts Love la vitesse et on secoue les boîtes, ça va mettre les voiles si ça pue la friteuse Conduite du genou, ça ressert en ligne droite Attends le sang Vas-y j''arrive Attends le sang Attends un peu Donnes moi une chaise Vas-y ramène J''ai pas besoin d''être debout, pour que je les baise Elle me dit T''es fou, tu veux me casser sans la bague sans la bague Marbella ou Braga, j''fais des vacances j''remonte du taga d''la bonne selha Hausse son Paul Shark, j''la fais tchouper dans la Bentayga cette petite kahba T''sais qu''pour 7500 j''fais les papiers à Air Harraga ? elle a trop la dalle J''attends rien de personne, j''fais des bails Des fois j''ai pleuré seul, sous bouteille Elle me fait la fille d''Insta, sous Lidl Elle se met fraîche dans sa chambre, en poubelle Elle me dit T''es fou, tu veux me casser sans la bague sans la bague Marbella ou Braga, j''fais des vacances j''remonte du taga d''la bonne selha Tu l''verras pas qu''on souffre Pas à n''importe qui on s''ouvre Dans le dur, dans le vrai, j''suis en panne de temps Que des coups de pfff, à tout je m''attends J''m''aère, je ère, y''a R Signé Moumzer d''la Zer Chabrand collée sur le ventre, j''apprends, à tout j''me relève Les tirs, des poteaux rentrants, j''me tire, ma petite m''attend J''perds le contrôle, t''façon j''m''enfous moi j''crois qu''j''suis gris J''tourne dans la zone, énervé le BM je l''ai plié Cali j''consomme, sur ma vie sinon j''vais vrillé Rien qu''ils chantonnent, ils ont parlés quand j''avais nié Dis moi où tu étais quand y''avait walou, nada Quand y''avait rien dans la CB, que je mettais des revers comme Nadal On barode avec le T, j''ai fini plus loin que Malaga J''ai tout mis dans l''OCB, un boulot j''te le fais harraga I''m not a hip?, I''m in a Bentley Gimme a tap now, can even end trip ? Whenever I tryna go, I''m around here I keep a blade like I''m Snipes Wesley You know where I come from, don''t try me They get Pop''d like Smoke or IP? Men, I can''t even stab his body They wont play tous les mecs de ma ville J''ai pris sur moi pour pas rapper en anglais, but I can''t held myself Drop if you really tryna spend this bag, hop in the jet or the brand-new ? Call all my niggas like Wesh le gang, they ain''t gonna play tough when I stop his bands They tryna eat all out of my plate, I tryna fuck all my hit J''veux péter le champ'' comme pilote de la F1, dans l''appart ça détaille fort de la S1 Dans le sac y''a le prix du Audi S1, on fait les cartes, on veut le code et on décaisse hein On est pas au flambe?, c''est la cité dans le 7R, on plait aux femmes, c''est la cité dans le secteur C''est la guerre, on a des traqueurs indétectables, on fait la paire, moi et mes couilles c''est le spectacle J''suis cramé comme Metah, chromé comme métal, la me-ar est létale, on bang on est fi dar On te rend méconnaissable, et faut qu''tu négocies, on veut remplir des salles, mes collègues aussi La j''suis d''humeur monotone, les schmitts qui passent, Molotov Que des colleurs, que des pocs, j''te croise, heynen, coup d''épaule Spliff au bec, j''suis sur la E17, 4 harbis dans un Audi break Tard la nuit j''m''enfume au Sour Diesel, avec j''fais des boulettes sur le jean Diesel J''roule à contre-sens pour finir dans l''top Ten, Amnésia, Hollandaise comme Robben La qualité t''inquiètes c''est de la Top Shelf, j''sors tout d''mes couilles y''a la douane française J''voulais pas m''salir la ligne est d''jà franchie, j''ressors blanchis, sourire d''affranchis Cali Weed scellé dans la saisie, ? j''avais d''jà senti Le ? m''appelle, j''le visser dans l''quart d''heure, qualité pas chère qui t''laisse parterre QN, La Haine, fumette dans l''artère, j''suis la 13014ème balle dans l''chargeur En freestyle, Le J m''appelle, donc j''arrive haut les curs J''me méfie de ceux qui ont du sucre dans la bouche et du venin dans le cur Qu''est-ce que la vie de rêve ? Voilà une question à méditer Des mythos m''ont promis de l''aide, au final c''est Dieu qui m''a facilité Ce mois-ci ça a rafalé combien d''fois ? Bienvenue à Marseille Y''a des mecs qui te cloueront sur une croix, si tu veux jouer l''immortel Mets pas ton nez dans mes affaires, tu veux nous faire mais tu tiens pas la route Je leur mets la patate à Fédor, quand on est faits d''or, sache qu''on ne craint pas la rouille J''ai des foutues pensées, elles sont ancestrales, et souvent j''me demande quand tout ça cessera Le viseur est sur toi quand passe sur l''estrade, comme l''envers du décor il est assez crade Frérot on a souffert jusqu''à l''insomnie, quand tu vois que le visage devient démoniaque Ils ont parlés trop vite moi j''ai rien compris, elle est bonne on l''a parfume à l''ammoniaque La vie cette pute elle a brisée des rêves, moi j''étais qu''un gosse qui essayait de tourner la roue Tourner la roue mais le monde est cruel comme aller bien sans Mel?, il faut finir le round J''connais sa mère pas son itinéraire, un désaccord et y''a mutinerie Soldat de Marseille avec un oeil de verre, pour tout prendre il te faudra plus qu''une vie C''est 1.3 Marseille bébé, c''est inné on fait TB En C.P ça fait TP, p''tit du CV froisse ta CB Week-end Monté-Carlo, suite et cuistot privé Mariage ?, fils de pute j''vais pas m''en privé Et vas-y vole petit pélican, amorti, lucarne, Bellingham J''fais mes comptes jusqu''à Alicante, elles sont jalouses mais on est polygames En Liga, eux ils sont reléguables, fais pas le raciste à la Monégasque Trop fier c''est moi qui régale, salon, yacht, en mode Monégasque Baisse pas les bras, redouble d''efforts, dans les difficultés on a vu les plus forts Pas de hasards, pas d''ironies du sort, j''suis pas à l''abri, pour moi envoie d''la ce-for J''suis dans l''sud, j''prends pas le nord, pendant j''coulais un bronze j''ai rêvé de l''or On est discrets, évite de parler fort, n''oublie pas d''où tu viens, surtout dans le Grand Fort? Sort couvert c''est glacial dehors J''suis dans l''nord, on descend vers le port J''refais tout, j''repense à tout Chacun ses atouts, sur le route j''mets des accoups Tu croyais quoi ? On est plusieurs dans ma tête, mais ? c''est moi le chef J''ai fais du bien, on peut aussi te faire du mal, me demande pas pour combien On est niah?, mais si j''dois t''niquer ta mère, on l''fera pas sur du Take Un petit pas, et tout le monde fait marche arrière, un peu comme cet Augustin J''m''en bas les couilles si t''as pas la ref'', le savoir est une arme que je brandis On fait pas d''exceptions à la règle, tant que la moula sonne, bah j''fais la diff'' Calibré dans l''bas du bâtiment, j''me suis tout en haut du Sofitel Affolant, j''ai fais et refais le tour de la ville, j''ai claqué des salaires dans des paires italiennes C''est banal ici on a tous bibi, gros zouave resert les loups, ceux qui ont soifs Un douze sur moi, faut pas qu''t''essayes, mec du CV quoi qu''il en soit Et même vacances, chemise en soie, ça reste à l''affût des poignets Ça vient de la Plaine jusqu''au Panier, gros c''est le 1.3 quoi qu''il en soit J''sais qu''toi t''es tout doux, nan fais pas celui qui aime les calibres On me dit Pour toi j''les fais dodo, mais tranquille, prends des années c''est pas la peine Merci le J pour l''assist'', toute la ville rep'' à l''appel J''tournais partout pour la SIM, minimum 2G pour l''appel 0.1 vers faire place, et laisse faire, on avait besoin d''espace, et d''espèces Donc on a pris des risques, pas à moi j''ai que les restes, la cité, je l''incarne, j''ai l''instinct Pour les faire, coups bas, là laisse moi un instant, j''ai la vista de Modri, Iniesta Tu connais j''suis dans les dièses, nan j''suis pas la pour le buzz, moi j''suis là que pour ma ville, mes ? dans l''cur Palpité ? Viens. Assieds toi, j''colle, manita reste à ta place, eux j''sais même pas t''es qui toi Faut que j''la quitte même si j''l''aime trop, tu sais que pour ma cité moi j''donne tout Ça presse la détente, shooter sur bécane, c''est plus comme avant ça tire dans l''tas Des minots chargés jusqu''aux dents, gantés La lumière viendra comme dans Billie Jean, Billie Billie T''as une vision de l''infini, 1.3.0 j''mets plus de limites La lumière viendra comme dans Billie Jean, Billie Billie Funambule, j''marche sur un fil, tu m''sors ton ? ça défile J''crois qu''t''as raison j''suis dans l''bail, et puis c''t''année on f''ra du sale RS4, j''mets les flammes, pour tous ces ? qui m''veulent du mal J''les shooterai tou-tou-tous, un par un je leur ferai du sale On fire, j''suis comme Jordan, j''suis dans les airs, rien que je fly Gambino, les z''hommes, www.vaniquertesmorts.com Ce soir elle s''est fait coquette, elle consomme ces morts j''ai niqué un vert pomme Ohh single de dia-dia, tu guettes ça ma baby J''aime pas les putas, c''est le 1.3, ça touche le RSA, ça roule en RS3 Ça fait tourner le four et ça pointe à la SPIP, cagoulé comme Hitch, tu t''prends un coup d''cric Y''a des bons et des putas et des bandeurs de hits, la vie ça va vite, parce que avant l''Covid Le jobbeur à décapsule, il fait tourner le shit J''suis au stade popopopopo Le virage m''a appellé, popopopopo, on fait trembler tout le Vél'' J''aurais allumé un pét'' avec si j''avais porté la flamme comme Jul J''crois que mehlich on va faire avec, si t''as le zemeh et pas d''bezoul Faut la Sacem à Cheb Khaled, ennemi public, Saddam, Ben Laden On s''en balek des strasses, des paillettes, et si ça flop j''décharge des palettes Quiquonque fait le King-Kong, et les balles du AK font du ping-pong P''tit con te jettes pas, j''mets des p''tits ponts, sort le cocktail et le chiffon Si tu veut la selha, y''a c''qu''il faut, des blancs, des noirs et des bicots D''or et d''plata sur le tifo, chico Pas d''estime, ils se ventent pour briller J''suis au rendez-vous même si j''suis trop grillé En Palestine, où vont mes prières Garder la face surtout faut pas plier Viens voir comment on est doués Le Loupe débarque, allez courez Tu fais le fou, on va te retrouver Viens voir comment on est chauds, vient voir comment on est capables J''débarque même dans ton ghetto, j''déverse ma haine, ça fait pa-pa Marseille, GTA, ? , corniche, Clio 3, réserve, ça c''cale à Cor-Beach ? Que des 30 ?, cagoule dans la biche?, sur les pleins capuches, c''est la rue ma biche Toi, tu parles des gens, tu m''intéresses pas, tu parles d''oseille t''es au RSA 13''Organisé comme le MS13, ça te ? au bar, ça met le GS3? Coeur cassé comme Bambi, fais un bisou à ton ? chérie j''t''en prie J''suis comme ceux qui m''écoutent, on s''est compris Coeur cassé comme Bambi J''reviens gazé et tant pis J''suis gazé et tant pis, j''suis comme ceux qui m''écoutent, on s''est compris Y''a pas de ?, on s''est compris, y''a pas de Ohhhhh Hamo C''est La Scampia, pose ton cul à terre Ohhhhh Hamo ? On compte plus les décès, ni les perquis'', ici tout est biaisé, c''est hyper triste Les mères ont les bras baissées et les pères prient, car ça rêve de Féfé et de berlines Trop de blessés, beaucoup de guerres crimes, la ville est aggressée, réput'' ternie M''empêche pas de rêver, tout est permis, je veux la LDC sous des zerbi? Du feu dans mes poèmes, moi j''n''ai plus le temps, jamais trop tard pour frapper Supporter de l''O.M, noir et musulman, j''suis le cauchemar de Galtier Faut qu''on rappe, faut qu''on les choque, faut qu''on les traumatise, moi j''me rappelle à l''époque ils me faisaient pas la bise Celle-là c''est pour Gordech? parti beaucoup trop tôt, bruits d''armes la gamberge du ghetto Khapta, t''as façon de parler tout tes mimiques ça passe pas, en vrai niquer des mères je sais pas j''suis bon qu''à ça Des sous contre du shit avec un tour de passe-passe avec un tour de passe-passe C''est plus fort que moi, j''vide la bouteille J''suis agen tous ces morts, c''est mort j''vide la bouteille Hier encore ça tirait, tirait, c''est des traîtres, c''est des Tiro? 6 heures du mat'', j''suis déchiré chiré, moi j''ai esquivé les gyros gyros J''suis Marseille, j''suis le 1.3 comme la vente de Coca'', atmosphère électrique j''me place comme Big Poppa C''est 13''Organisé, DZ mon avocat, coup d''schlass dans la ceinture, j''te termine au Toka'' J''fais du rap de crapule, pas de nez dans la coco, j''arrive comme à l''ancienne, VR6 full Croco Sauvage comme Dothraki, barbare comme Khal Drogo, ça m''entoure des gros coups, mafieux comme la Mocro? J''arrive comme Marcheur Blanc, j''t''ouvre le crâne à l''épée, j''compte plus les heures de gardav'', ni les heures de TP J''suis le Réal de Madrid comme Kakà ou Pépé, violent depuis le CP, PGP full C.P On fuck la Fashion Week, pas de dégaine de pétasse, on braque après on s''casse, on brasse après on tasse Chez nous y''a pas de Féfé, R19 j''mets les gaz, ? foudroie la prod et Sysa sur Pégase 11.43 j''ai décapoté, ça bipe sur le bleu toute la journée J''ai vidé la ''teille j''suis assommé, brolique sous la ceinture j''deviens jnouné Dis moi où Dadinho n''a pas allumé, j''fais parti des négros miraculés J''en ai bégayé, on a assumé, ils m''appellent le sang, c''est des enculés Comme le J j''représente la zone, c''est seulement le frérot qui a fini en zone Garder le moral ça brise pas les hommes, 13''Organisé ça fait l''effet d''une bombe D''une bombe, d''une bombe, 13''Organisé ça fait l''effet d''une bombe On verbalise puis on pénalise, quand c''est la crise, que la moula qui nous canalise Pas tout seul j''suis en bande organisée, on est là pour traumatiser J''roule une gros joint sous la lune, j''m''évade même sans thunes, la tête dans les étoiles, les pieds sur le bitume Entre le bien et le mal, j''sais plus où je me situe, j''rêve de nouveaux horizons même si ici ça tue Et ouais ça tue en rafale, comme aux Favelas, matin midi soir ça remplit le sac de mapessa On porte nos couilles sur le terrain, comme Makélélé, les gars ça tire, même les condés dansent la Macarena J''sors du four ça sent le poulet, y''a les schmitts qui tournent, on s''en bat les-les Y''a les schmitts qui tournent, on s''en bat les-les, allumer le feu, c''que tu fais c''est laid C''que tu fais c''est pété, moi j''écoute pas, accepte les critiques et devient plus fort On va choquer du monde, la y''a toute la ciudad, j''ai mis les lunettes à Mad Max avec la ? En ce moment couci-couça, en manque de moula donc j''ai bipé le ''zin Pardonner les fous c''est un truc que j''fais pas, y''a pas de silencieux donc j''utilise un coussin Ah haut les mains haut les mains, t''as claqué tout ton biff maintenant tu veux le mien T''auras rien t''auras rien, y''a que quand tout va bien qu''t''es le frangin t''es le frangin Doucement, pousse toi et assis toi regarde faire les anciens J''ai les crocs, j''ai la rage comme un chien quand ça touche à mon pain Ya des frères qui sont tombés bêtement Jai pas changé depuis mes 20 ans Jarrive 1.3.9, je crie Vatos Cest le D est-ce que tu me remets maintenant ?? Jarrive inattendu comme un coup de tazer Jarrive tendu comme le cul de ta sur Ouais dans les rues de Marseille elle arrive par la poste, on a plus de passeur Regarde les bandits brasser Je béni lavenir, je maudis le passé Sur le prado je conduis gazé Je reviens de loin comme Jean-Louis Gasset Je suis un couche tard, jai ma tare-Gui Ne la touche pas, cest ma barbie Marseille freestyle légendaire dans la poche jai du méchant teh tu serres'
$ws.Cells.Item(6, 3).Value = 'On se bat pour ne pas perdre Sortis de notre réserve comme des Apaches On crée la psychose au milieu des visages pâles Chevauchée des indigènes des quartiers sales Gardiens de la dernière tribu des XXX Sentiment d''insécurité dans ces black nights On vient prendre le mic et on fait chauffer la dancehall Qu''la cavalerie débarque, c''st déjà trop tard Le son des Mohicans a fait des ravags J''viens d''là où on est tous aigris et les sourires s''monnayent Où l''taux d''alcool, shit, remplace donc le monoï Là où la BAC flique les canailles Aussi louche que George Bush en plein centre d''Hanoï Où l''on rêve d''être riche comme le baron d''Brunei Ça reste en stand-by XXX venu d''Dubaï J''viens d''là où on a tous soif d''idéal Mais trop épris des halls et l''amour est sous scellés Où la patrie a tourné l''dos à ses fistons Et c''est love story entre bad life et béton J''arrive solide comme my man lors d''une XXX Inonde la street de ma zik comme XXX demande à Saher Tous besoin d''air, mes mans ont l''goût du risque La rage se lit dans nos yeux et s''écoute sur disque Mon son c''est Bagdad, Boss One au mic c''est vendetta Préviens l''État qu''c''est Chroniques XXX ça Sortis de notre réserve comme des Apaches On crée la psychose au milieu des visages pâles Chevauchée des indigènes des quartiers sales Gardiens de la dernière tribu des XXX Sentiment d''insécurité dans ces black nights On vient prendre le mic et on fait chauffer la dancehall Qu''la cavalerie débarque, c''est déjà trop tard Le son des Mohicans a fait des ravages C''est nous ton western, les tournages en plein Far West Shoota tous les James West au Nord au Sud, à l''Est, Ouest Demande à Géronimo si on plaisante, mec Demande à Blueberry qu''il fasse tourner son pèt'', merde Kilomètre cinq sept, les XXX s''battent en duel On dégaine son revolver, dans la rue ça fait des bang-bang Cow-boy d''la BAC en civil, faut pas qu''tu testes, ah Remballe ton bleu gyrophare et ne traîne pas dans nos ruelles Quand le XXX est lâché, XXX donc sur la Canebière Charge de CRS matraquez-moi ces métèques La fumée, odeur de bomba lacrymogène On avance les yeux fermés On se bat pour ne pas perdre Sortis de notre réserve comme des Apaches On crée la psychose au milieu des visages pâles Chevauchée des indigènes des quartiers sales Gardiens de la dernière tribu des XXX Sentiment d''insécurité dans ces black nights On vient prendre le mic et on fait chauffer la dancehall Qu''la cavalerie débarque, c''est déjà trop tard Le son des Mohicans a fait des ravages Tu sais, c''est sur la route de la faiblesse que j''m''ennuie Enfant des quartiers, là où on s''bat pour forcer l''avenir Là où se dealent la merde et la tristesse Là où les p''tites tafioles se font prendre de vitesse Là où se fume la pharmacie des quartiers Là où Clara Morgane est le rêve des cages d''escaliers Tu sais, j''ai vu danser la haine Chez nous ils attendent que Jacques Chirac purge sa peine De là où la liberté arrache des sacs L''espoir se joue aux cartes et la galère échoue sa barque Je viens d''où? Je viens de ces coins d''France Où le Zé Péquenho connaît l''allure des baraques Ils veulent tous porter la coupe du ghetto La sirène des schmitts devient la symphonie des khos Là où les bambins veulent être des boss B. One et Toko Blaze, c''est signé El Sarazino Sortis de notre réserve comme des Apaches On crée la psychose au milieu des visages pâles Chevauchée des indigènes des quartiers sales Gardiens de la dernière tribu des XXX Sentiment d''insécurité dans ces black nights On vient prendre le mic et on fait chauffer la dancehall Qu''la cavalerie débarque, c''est déjà trop tard Le son des Mohicans a fait des ravages On se bat pour ne pas perdre Toko Blaze On se bat pour ne pas perdre La Swija On se bat pour ne pas perdre Chroniques de Mars On se bat pour ne pas perdre'
$ws.Cells.Item(7, 3).Value = 'en espagnol Mental d''acier comme un vrai champion Depuis minot on rêve d''être un champion Quitter l''quartier, devenir champion Champion, champion Gravir les marches jusqu''au sommet champion Barrages de titres et réussir champion Dur parcours de l''enfant devenu champion Champion, champion On est tous sortis du ventre de notre mère Dans le brouillard, incertitude, on cherche nos repères Chacun choisit sa voie, chacun son itinéraire Sans XXX on s''doit de mieux faire que nos pèrs Traverser les océans, travrser les mers Franchir le froid, le vent violent et le désert Marcher sur les routes semées d''embûches, de galères Le talent passe toujours au-delà des frontières Dix, dix, dix On porte tous le numéro dix Sur l''terrain dix, dix, dix On voit que des numéros dix Dans l''ghetto dix, dix, dix Nous n''sommes que des numéros dix Porte le numéro dix On a la rage de réussir Dix, dix, dix Gagnant comme un numéro dix Sur l''terrain dix, dix, dix Dans le dos le numéro dix Dans l''ghetto dix, dix, dix Magie de l''enfant prodige Légendaire numéro dix Porte le dix comme un artiste J''ai des problèmes psychologiques d''un loco À la tranquille, je pose le seize sur le morceau de Toko Ils finissent mal, un peu comme Patrice Loko Ce monde est fou, il faut s''unir et faire ça entre potos Tu n''obtiens rien sans faire de sacrifices Enfants d''la rue, d''la colère un peu comme Mary Pierce On n''avait rien ici, l''bonheur c''était une rue, un ballon Bouteille de Québec dans l''alim, on volait tous les bonbons Stade Vélodrome, pesage ou populaire Traverse la ville avec les frères comme on traverse la Terre On s''est perdus sans retrouver les chemin Fais pas comme nous, comme le Suprême tu es le monde de demain Maillot floqué, tu as l''numéro dix Le chef d''orchestre de ta vie mais c''est toi qui décides Personne pour le faire à ta place Évite les pièges de cette vie ou finis broyé dans la masse Dix, dix, dix On porte tous le numéro dix Sur l''terrain dix, dix, dix On voit que des numéros dix Dans l''ghetto dix, dix, dix Nous n''sommes que des numéros dix Porte le numéro dix On a la rage de réussir Dix, dix, dix Gagnant comme un numéro dix Sur l''terrain dix, dix, dix Dans le dos le numéro dix Dans l''ghetto dix, dix, dix Magie de l''enfant prodige Légendaire numéro dix Porte le dix comme un artiste Livrés à nous-mêmes comme des bonhommes On s''est battus pour exaucer nos rêves de gosses Contourner les remparts ou enfoncer des portes On est allés l''chercher, pas attendu qu''on nous l''donne Est-ce que tu seras prêt si ton heure sonne Le bon moment, le bon endroit, la bonne personne Médaillé d''or ou diplômé, remonte le score Être ambitieux sans marcher sur la tête des autres Entre les barbus, les FAF Les têtes capuchées, l''crack Sois Maradonna à Naples, adulé dans les ruelles crades J''porte le même prénom que Messi en vrai Le kiss dans l''équipe, non juste un con d''plus dont la vie a fait Ranger son éthique, au cou la marque de la dague Beaucoup d''palabres et d''arnaques J''avais tout, dire qu''à mon doigt s''estompe la la marque de la bague Petit, tu rates une passe, deux passes, soudain l''stade se vite Encore plus si c''est penalty Les potos t''raviront le titre On mène le jeu quand leur État s''attache à gruger l''prolo Iconoclaste par défaut, demande au Muge et Toko Même en disgrâce, porte le dix gars Et kicke ça même si l''quidam jette le vice et canarde ton équipe Sache mec, que l''terrain est miné d''avance Comme les boss de la FIFA, ils ne visent que bitch et argent C''est pour ceux qui portent le dix au cur, ouais ceux d''ma ville Pour tout c''qui n''a pas de prix comme les yeux de ma fille Dix, dix, dix On porte tous le numéro dix Sur l''terrain dix, dix, dix On voit que des numéros dix Dans l''ghetto dix, dix, dix Nous n''sommes que des numéros dix Porte le numéro dix On a la rage de réussir Dix, dix, dix Gagnant comme un numéro dix Sur l''terrain dix, dix, dix Dans le dos le numéro dix Dans l''ghetto dix, dix, dix Magie de l''enfant prodige Légendaire numéro dix Porte le dix comme un artiste Mental d''acier comme un vrai champion Depuis minot on rêve d''être un champion Quitter l''quartier, devenir champion Champion, champion Gravir les marches jusqu''au sommet champion Barrages de titres et réussir champion Dur parcours de l''enfant devenu champion'
$ws.Cells.Item(8, 3).Value = 'Origines Maintien du lien avec nos racines Éloignées comme un parent Aller chercher la réponse, combler ce vide Héritier d''une double culture Aussi riche que l''enfant métis Sur la terre de tes ancêtres démarre ton histoire Ici, tes nouveaux liens s''tissent Mes origines, tes origines, leurs origines wey wey wey, wey Nos origines, vos origines, ses origines wey, wey, wey Les origines, seules origines, quelles origines wey wey wey, wey Comme origines, pour origines, toutes origines wey, wey, wey Chaud et humide Comme le climat sous les tropiques Saison des pluies, les pages se noircissent Toutes ces muses m''inspirent Tension, réflechir La démarche reste artistique Griot rebelle d''une tradition orale, musicale et poétique Vibration exotique Afro-style, la mode est cyclique Expatrié, connexion de la diaspora dans la city Affronter les préjugés Combattre le racisme Resserrer les liens, se serrer les coudes, éviter l''communautarisme Mes origines, tes origines, leurs origines wey wey wey, wey Nos origines, vos origines, ses origines wey, wey, wey Les origines, seules origines, quelles origines wey wey wey, wey Comme origines, pour origines, toutes origines wey, wey, wey À Marseille c''est, yeah bana bana Dans toute l''Europe, yeah bana bana Dans toute l''Afrique, yeah bana bana États-Unis, yeah bana bana Amérique Latine, yeah bana bana Dans les Caraïbes, yeah bana bana Océan Indien, yeah bana bana Océanie, yeah bana bana Dans toute l''Asie, yeah bana bana Royaume-Uni, yeah bana bana À Abidjan, yeah bana bana À Yaoundé, yeah bana bana À Conakry, yeah bana bana À Niamey, yeah bana bana À Librevile, yeah bana bana Tout l''monde ici, yeah bana bana Mes origines, tes origines, leurs origines wey wey wey, wey Nos origines, vos origines, ses origines wey, wey, wey Les origines, seules origines, quelles origines wey wey wey, wey Comme origines, pour origines, toutes origines wey, wey, wey'
$ws.Cells.Item(9, 3).Value = 'Eh, Tok'' Blaze XXX Pouss Pouss Laisse-moi passer Le coach a décidé ce soir de ne pas m''aligner Je m''suis entraîné dur pour que mon show soit carré Toko! Remplaçant de luxe, toujours pas titularisé Les crampons sont cirés, le microphone prêt à fumer Faut rester concentré, ne jamais se décourager Malgré la concurrence, on doit toujours persévérer Sur l''terrain, sur la scène, l''objectif ma place à gagner J''viens pas chauffer le banc, moi aussi j''veux participer Pouss Pouss Pas l''temps de douter Pouss Pouss Pas l''temps d freiner Pouss Pouss La devis c''est Droit au but, faut me laisser passer Pouss Pouss Pas l''temps de reculer Pouss Pouss De se retourner Pouss Pouss Reste là sur le bord de la touche et laisse-moi passer Pouss Pouss Laisse-moi passer Pouss Pouss Laisse-moi passer Pouss Pouss Laisse-moi passer Pouss Pouss Laisse-moi passer Si t''es en baisse de régime, suspendu ou blessé Si t''as pas convaincu je vais devoir te remplacer à moi, à moi T''as perdu ta place, une chance à ne pas gâcher XXX les changements Maintenant, à mon tour de briller D''Abidjan à Marseille, de Cotonou à Yaoundé Ballon d''or franco-africain pour notre réalité On tue l''match dans les charts, aucune défense ne peut stopper Sur la piste de danse, c''est la victoire assurée Pouss Pouss Pas l''temps de douter Pouss Pouss Pas l''temps de freiner Pouss Pouss La devise c''est Droit au but, faut me laisser passer Pouss Pouss Pas l''temps de reculer Pouss Pouss De se retourner Pouss Pouss Reste là sur le bord de la touche et laisse-moi passer Pouss Pouss Laissez-moi passer Pouss Pouss Laissez-moi passer Pouss Pouss Laissez-moi passer Pouss Pouss Laissez-moi passer Pouss Pouss Pouss Pouss Pouss Pouss Pouss Pouss XXX Pouss Pouss Pas l''temps de douter Pouss Pouss Pas l''temps de freiner Pouss Pouss La devise c''est Droit au but, faut me laisser passer Pouss Pouss Pas l''temps de reculer Pouss Pouss De se retourner Pouss Pouss Reste là sur le bord de la touche et laisse-moi passer XXX'
$ws.Cells.Item(10, 3).Value = 'Condamné pour un délit de sale gueule Célibataire endurci dans le deux pièces dun immeuble Enfance de fils unique, grandir sans frère ni sur Père succomba, Mère ne supportait plus dêtre veuve Beau père alcoolique, placement en famille daccueil Dans la cour de récré, toujours de ceux qui restaient seuls Sans chagrin damour, vibration dun premier flirt Goût du baiser, odeur, parfum dune amoureuse On se sent seul, dans le même cas on est si nombreux Isolés les uns des autres, on communique peu Eux quand ils rentraient chez eux, le regard amoureux Tout pouvait sécrouler autour deux On sent seul en même temps de plus en plus nombreux Pas besoin de grand chose pour se sentir heureux Mais on se fait déjà trop vieux a beau fermer les yeux Espérant quun jour tout ira mieux Comme ces instants où jrêvais dvoir la ville en feu Après altercation avec tous ces hommes en bleu Question drespect, on mettait not parole en jeu Brisant les clichés de gars venus dun quartier dangereux Quand ils rentraient chez eux, jtrônais comme un roi sur mes bancs A forcer lentrée de tant de soirées dansantes Le hip-hop ouais cest ça quon avait dans lsang Autour dmoi, peu dlettres, mon gars ça parlait en francs Cette impression dvoir Mars séparé en deux Ceux qui tirent les fils et ceux qui vivent sans enjeu Qui voient jamais rien dneuf et se niquent entre eux Portant sur les épaules, un gros passif de contentieux Jrevois lvillage où jadis jvivais enfant Puis lescalade des noises et son tracé en pente Ma place dans le quartier, assailli par les cancans Tes lèvres sur mon front, et mes blessures qui guérissent dans ltemps On se sent seul, dans le même cas on est si nombreux Isolés les uns des autres, on communique peu Eux quand ils rentraient chez eux, le regard amoureux Tout pouvait sécrouler autour deux On sent seul en même temps de plus en plus nombreux Pas besoin de grand chose pour se sentir heureux Mais on se fait déjà trop vieux a beau fermer les yeux Espérant quun jour tout ira mieux Victime anonyme à la recherche du bonheur Dans lobscurité, à se dire que la solitude fait peur Lindividualité te condamne à vivre seul Te donne limpression dêtre passé à deux doigts du bonheur Dans lindifférence, je prends place dans lascenseur Quotidien, se coucher tard pour se lever de bonne heure Souvrir une boite, sendormir devant le téléviseur Personne ne se souciera si un jour je pars avant lheure On se sent seul, dans le même cas on est si nombreux Isolés les uns des autres, on communique peu Eux quand ils rentraient chez eux, le regard amoureux Tout pouvait sécrouler autour deux On sent seul en même temps de plus en plus nombreux Pas besoin de grand chose pour se sentir heureux Mais on se fait déjà trop vieux a beau fermer les yeux Espérant quun jour tout ira mieux'
$ws.Cells.Item(11, 3).Value = 'Aigle appelle Caïman, Aigle appelle Caïman C''est Radio Bananier Sport, en direct du stade de la mission Pour un grand match au sommet entre deux équipes chocs Les Planteurs de Manioc contre l''Olympique Mandingue Les effectifs sont armés, la partie s''annonce dangereuse Balle à terre, les Olympiens balancent L''OM lance une offensive La défens reste complètemnt figée, la surface de réparation est claire Transversale chirurgicale Dans les dix-huit mètres, c''est danger Incendie oh, appelez les pompiers, le gardien de but est paralysé Le ballon file Ce n''est pas possible, le ballon a quitté le terrain pour aller se loger dans l''enceinte même de l''église Sachant le problème qu''il existe entre l''église et l''hôtel de ville, il va être difficile pour nous de récupérer le ballon Tout à fait Madou, et Monsieur le curé a pris en otage le ballon et réclame en échange de celui-même Une bouteille de gin et un plat de riz de chez Mama Kounté Le temps que les hautes autorités du football se soumettent à ces revendications Nous allons vous rendre l''antenne pour une page de publicité C''est Madou Diallo et Rigobert Koulé en direct du stade de la mission À vous les studios XXX À tout à l''heure'
